$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E ("estado"), styled like the rest of the header row
$ws.Range("E1").Value = "estado"
$ws.Range("E1").NumberFormat = "@"

# Existing rows 2 and 3 get a new "estado" value in column E (default/general format)
$ws.Range("E2").Value = "Creado"
$ws.Range("E3").Value = "Creado"

# New row 4 data (text-formatted, like the other data rows)
$ws.Range("A4").Value = "tipobusqueda3"
$ws.Range("B4").Value = "numeroDocumento3"
$ws.Range("C4").Value = "garantia3"
$ws.Range("D4").Value = "motivo3"
$ws.Range("E4").Value = "Creado"
$ws.Range("A4:E4").NumberFormat = "@"

# New row 5 data (text-formatted in A:D, default/general format in E)
$ws.Range("A5").Value = "tipobusqueda4"
$ws.Range("B5").Value = "numeroDocumento4"
$ws.Range("C5").Value = "garantia4"
$ws.Range("D5").Value = "motivo4"
$ws.Range("E5").Value = "Creado"
$ws.Range("A5:D5").NumberFormat = "@"

# Selection moves to E4, matching the recorded state
$ws.Range("E4").Select()
